# The "Flight Times" sheet holds a symmetric city-to-city flight-time matrix
# (rows 2-8 / cols B-H). Placeholder cells containing -1 are filled in with
# the mirrored (transposed) time value from the opposite side of the
# diagonal, and formatted the same way the already-populated time cells are
# (the "h:mm" time number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3" = 0.17361111111111113
    "E3" = 0.18055555555555555
    "F3" = 0.22916666666666666
    "G3" = 0.27083333333333331
    "H3" = 0.17361111111111113

    "C4" = 0.17361111111111113
    "F4" = 0.22222222222222221
    "G4" = 0.15625
    "H4" = 0.11458333333333333

    "C5" = 0.18055555555555555
    "G5" = 0.072916666666666671

    "C6" = 0.22916666666666666
    "D6" = 0.22222222222222221
    "G6" = 0.16319444444444445

    "C7" = 0.27083333333333331
    "D7" = 0.15625
    "E7" = 0.072916666666666671
    "F7" = 0.16319444444444445

    "C8" = 0.17361111111111113
    "D8" = 0.11458333333333333
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    $cell.NumberFormat = "h:mm"
}
